# Adds the new "payment_address" test-data sheet, refreshes the
# register_test selection/active-tab state, and leaves the new sheet
# active with its own selection - matching the upstream commit
# "updating test data for upcoming tests".

$wb = $excel.ActiveWorkbook

# --- register_test: move the remembered selection -------------------------
$regSheet = $wb.Worksheets.Item("register_test")
[void]$regSheet.Range("F21").Select()

# --- add the new sheet after the last existing one -------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "payment_address"

# --- header row --------------------------------------------------------
$newSheet.Range("A1").Value = "first_name"
$newSheet.Range("B1").Value = "last_name"
$newSheet.Range("C1").Value = "address"
$newSheet.Range("D1").Value = "city"
$newSheet.Range("E1").Value = "post_code"
$newSheet.Range("F1").Value = "country"
$newSheet.Range("G1").Value = "region"

# --- sample data row -----------------------------------------------------
$newSheet.Range("A2").Value = "test1"
$newSheet.Range("B2").Value = "test"
$newSheet.Range("C2").Value = "testing address1"
$newSheet.Range("D2").Value = "test city"
$newSheet.Range("E2").Value = 54645
$newSheet.Range("F2").Value = "India"
$newSheet.Range("G2").Value = "Karnataka"

# --- reuse register_test's header/data formatting (yellow fill + border +
#     centered for the header, centered fill-less for the data row) --------
$regSheet.Range("A1:B1").Copy()
$newSheet.Range("A1:G1").PasteSpecial(-4122)
$regSheet.Range("A2:B2").Copy()
$newSheet.Range("A2:G2").PasteSpecial(-4122)

# --- widen the "address" column to fit its content -------------------------
$newSheet.Columns.Item(3).ColumnWidth = 13.25

# --- leave the new sheet active with its own remembered selection ----------
[void]$newSheet.Range("F10").Select()
